$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete fg_mf_total column (AB) - shrinks dimension to column AA
$ws.Range("AB1:AB22").Delete() | Out-Null

# --- Header row ---
$ws.Cells.Item(1, 1).Value = "comp_name"
$ws.Cells.Item(1, 2).Value = "underiv_comp_name"
$ws.Cells.Item(1, 3).Value = "iupac_name"
$ws.Cells.Item(1, 4).Value = "molecular_formula"
$ws.Cells.Item(1, 5).Value = "canonical_smiles"
$ws.Cells.Item(1, 6).Value = "molecular_weight"
$ws.Cells.Item(1, 7).Value = "xlogp"
$ws.Cells.Item(1, 8).Value = "el_C"
$ws.Cells.Item(1, 9).Value = "el_H"
$ws.Cells.Item(1, 10).Value = "el_O"
$ws.Cells.Item(1, 11).Value = "el_mf_C"
$ws.Cells.Item(1, 12).Value = "el_mf_H"
$ws.Cells.Item(1, 13).Value = "el_mf_O"
$ws.Cells.Item(1, 14).Value = "fg_C-aliph"
$ws.Cells.Item(1, 15).Value = "fg_C-arom"
$ws.Cells.Item(1, 16).Value = "fg_alcohol"
$ws.Cells.Item(1, 17).Value = "fg_carboxyl"
$ws.Cells.Item(1, 18).Value = "fg_ester"
$ws.Cells.Item(1, 19).Value = "fg_ether"
$ws.Cells.Item(1, 20).Value = "fg_ketone"
$ws.Cells.Item(1, 21).Value = "fg_mf_C-aliph"
$ws.Cells.Item(1, 22).Value = "fg_mf_C-arom"
$ws.Cells.Item(1, 23).Value = "fg_mf_alcohol"
$ws.Cells.Item(1, 24).Value = "fg_mf_carboxyl"
$ws.Cells.Item(1, 25).Value = "fg_mf_ester"
$ws.Cells.Item(1, 26).Value = "fg_mf_ether"
$ws.Cells.Item(1, 27).Value = "fg_mf_ketone"

# --- Data rows (reordered + columns reshuffled) ---
# Row 2
$ws.Cells.Item(2, 1).Value = "stearic acid, tms derivative"
$ws.Cells.Item(2, 2).Value = "stearic acid"
$ws.Cells.Item(2, 3).Value = "octadecanoic acid"
$ws.Cells.Item(2, 4).Value = "C18H36O2"
$ws.Cells.Item(2, 5).Value = "CCCCCCCCCCCCCCCCCC(=O)O"
$ws.Cells.Item(2, 6).Value = 284.5
$ws.Cells.Item(2, 7).Value = 7.4
$ws.Cells.Item(2, 8).Value = 18
$ws.Cells.Item(2, 9).Value = 36
$ws.Cells.Item(2, 10).Value = 2
$ws.Cells.Item(2, 11).Value = 0.7599226713532512
$ws.Cells.Item(2, 12).Value = 0.1275500878734622
$ws.Cells.Item(2, 13).Value = 0.1124710017574693
$ws.Cells.Item(2, 14).Value = 17
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 0.8417117750439367
$ws.Cells.Item(2, 22).Value = 0
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0.1582319859402461
$ws.Cells.Item(2, 25).Value = 0
$ws.Cells.Item(2, 26).Value = 0
$ws.Cells.Item(2, 27).Value = 0

# Row 3
$ws.Cells.Item(3, 1).Value = "benzene-1,2-diol, deriv"
$ws.Cells.Item(3, 2).Value = "benzene-1,2-diol"
$ws.Cells.Item(3, 3).Value = "benzene-1,2-diol"
$ws.Cells.Item(3, 4).Value = "C6H6O2"
$ws.Cells.Item(3, 5).Value = "C1=CC=C(C(=C1)O)O"
$ws.Cells.Item(3, 6).Value = 110.11
$ws.Cells.Item(3, 7).Value = 0.9
$ws.Cells.Item(3, 8).Value = 6
$ws.Cells.Item(3, 9).Value = 6
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(3, 11).Value = 0.6544909635818728
$ws.Cells.Item(3, 12).Value = 0.05492689129052766
$ws.Cells.Item(3, 13).Value = 0.2906003087821269
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 6
$ws.Cells.Item(3, 16).Value = 2
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 0
$ws.Cells.Item(3, 22).Value = 0.6911088911088911
$ws.Cells.Item(3, 23).Value = 0.3089092725456362
$ws.Cells.Item(3, 24).Value = 0
$ws.Cells.Item(3, 25).Value = 0
$ws.Cells.Item(3, 26).Value = 0
$ws.Cells.Item(3, 27).Value = 0

# Row 4
$ws.Cells.Item(4, 1).Value = "1-monopalmitin, 2tms derivative"
$ws.Cells.Item(4, 2).Value = "1-monopalmitin"
$ws.Cells.Item(4, 3).Value = "2,3-dihydroxypropyl hexadecanoate"
$ws.Cells.Item(4, 4).Value = "C19H38O4"
$ws.Cells.Item(4, 5).Value = "CCCCCCCCCCCCCCCC(=O)OCC(CO)O"
$ws.Cells.Item(4, 6).Value = 330.5
$ws.Cells.Item(4, 7).Value = 6.3
$ws.Cells.Item(4, 8).Value = 19
$ws.Cells.Item(4, 9).Value = 38
$ws.Cells.Item(4, 10).Value = 4
$ws.Cells.Item(4, 11).Value = 0.6904962178517398
$ws.Cells.Item(4, 12).Value = 0.1158971255673222
$ws.Cells.Item(4, 13).Value = 0.1936338880484115
$ws.Cells.Item(4, 14).Value = 17
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 2
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 1
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0.7215098335854766
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0.1029167927382754
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0.1756006051437216
$ws.Cells.Item(4, 26).Value = 0
$ws.Cells.Item(4, 27).Value = 0

# Row 5
$ws.Cells.Item(5, 1).Value = "1-monooleoylglycerol, 2tms derivative"
$ws.Cells.Item(5, 2).Value = "1-monooleoylglycerol"
$ws.Cells.Item(5, 3).Value = "2,3-dihydroxypropyl (z)-octadec-9-enoate"
$ws.Cells.Item(5, 4).Value = "C21H40O4"
$ws.Cells.Item(5, 5).Value = "CCCCCCCCC=CCCCCCCCC(=O)OCC(CO)O"
$ws.Cells.Item(5, 6).Value = 356.5
$ws.Cells.Item(5, 7).Value = 6.5
$ws.Cells.Item(5, 8).Value = 21
$ws.Cells.Item(5, 9).Value = 40
$ws.Cells.Item(5, 10).Value = 4
$ws.Cells.Item(5, 11).Value = 0.7075203366058906
$ws.Cells.Item(5, 12).Value = 0.1130995792426367
$ws.Cells.Item(5, 13).Value = 0.1795119214586255
$ws.Cells.Item(5, 14).Value = 19
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 2
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 18).Value = 1
$ws.Cells.Item(5, 19).Value = 0
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 0.7419270687237027
$ws.Cells.Item(5, 22).Value = 0
$ws.Cells.Item(5, 23).Value = 0.09541093969144461
$ws.Cells.Item(5, 24).Value = 0
$ws.Cells.Item(5, 25).Value = 0.1627938288920056
$ws.Cells.Item(5, 26).Value = 0
$ws.Cells.Item(5, 27).Value = 0

# Row 6
$ws.Cells.Item(6, 1).Value = "9-octadecenoic acid, (e)-, deriv"
$ws.Cells.Item(6, 2).Value = "9-octadecenoic acid, (e)-"
$ws.Cells.Item(6, 3).Value = "(e)-octadec-9-enoic acid"
$ws.Cells.Item(6, 4).Value = "C18H34O2"
$ws.Cells.Item(6, 5).Value = "CCCCCCCCC=CCCCCCCCC(=O)O"
$ws.Cells.Item(6, 6).Value = 282.5
$ws.Cells.Item(6, 7).Value = 6.5
$ws.Cells.Item(6, 8).Value = 18
$ws.Cells.Item(6, 9).Value = 34
$ws.Cells.Item(6, 10).Value = 2
$ws.Cells.Item(6, 11).Value = 0.7653026548672566
$ws.Cells.Item(6, 12).Value = 0.121316814159292
$ws.Cells.Item(6, 13).Value = 0.1132672566371681
$ws.Cells.Item(6, 14).Value = 17
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 0
$ws.Cells.Item(6, 20).Value = 0
$ws.Cells.Item(6, 21).Value = 0.8405345132743363
$ws.Cells.Item(6, 22).Value = 0
$ws.Cells.Item(6, 23).Value = 0
$ws.Cells.Item(6, 24).Value = 0.1593522123893805
$ws.Cells.Item(6, 25).Value = 0
$ws.Cells.Item(6, 26).Value = 0
$ws.Cells.Item(6, 27).Value = 0

# Row 7
$ws.Cells.Item(7, 1).Value = "octadecanoic acid, deriv"
$ws.Cells.Item(7, 2).Value = "octadecanoic acid"
$ws.Cells.Item(7, 3).Value = "octadecanoic acid"
$ws.Cells.Item(7, 4).Value = "C18H36O2"
$ws.Cells.Item(7, 5).Value = "CCCCCCCCCCCCCCCCCC(=O)O"
$ws.Cells.Item(7, 6).Value = 284.5
$ws.Cells.Item(7, 7).Value = 7.4
$ws.Cells.Item(7, 8).Value = 18
$ws.Cells.Item(7, 9).Value = 36
$ws.Cells.Item(7, 10).Value = 2
$ws.Cells.Item(7, 11).Value = 0.7599226713532512
$ws.Cells.Item(7, 12).Value = 0.1275500878734622
$ws.Cells.Item(7, 13).Value = 0.1124710017574693
$ws.Cells.Item(7, 14).Value = 17
$ws.Cells.Item(7, 15).Value = 0
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = 0
$ws.Cells.Item(7, 19).Value = 0
$ws.Cells.Item(7, 20).Value = 0
$ws.Cells.Item(7, 21).Value = 0.8417117750439367
$ws.Cells.Item(7, 22).Value = 0
$ws.Cells.Item(7, 23).Value = 0
$ws.Cells.Item(7, 24).Value = 0.1582319859402461
$ws.Cells.Item(7, 25).Value = 0
$ws.Cells.Item(7, 26).Value = 0
$ws.Cells.Item(7, 27).Value = 0

# Row 8
$ws.Cells.Item(8, 1).Value = "9,12-octadecadienoic acid (z,z)-, tms derivative"
$ws.Cells.Item(8, 2).Value = "9,12-octadecadienoic acid (z,z)-"
$ws.Cells.Item(8, 3).Value = "(9z,12z)-octadeca-9,12-dienoic acid"
$ws.Cells.Item(8, 4).Value = "C18H32O2"
$ws.Cells.Item(8, 5).Value = "CCCCCC=CCC=CCCCCCCCC(=O)O"
$ws.Cells.Item(8, 6).Value = 280.4
$ws.Cells.Item(8, 7).Value = 6.8
$ws.Cells.Item(8, 8).Value = 18
$ws.Cells.Item(8, 9).Value = 32
$ws.Cells.Item(8, 10).Value = 2
$ws.Cells.Item(8, 11).Value = 0.7710342368045648
$ws.Cells.Item(8, 12).Value = 0.1150356633380885
$ws.Cells.Item(8, 13).Value = 0.1141155492154066
$ws.Cells.Item(8, 14).Value = 17
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = 0
$ws.Cells.Item(8, 19).Value = 0
$ws.Cells.Item(8, 20).Value = 0
$ws.Cells.Item(8, 21).Value = 0.8396398002853066
$ws.Cells.Item(8, 22).Value = 0
$ws.Cells.Item(8, 23).Value = 0
$ws.Cells.Item(8, 24).Value = 0.1605456490727532
$ws.Cells.Item(8, 25).Value = 0
$ws.Cells.Item(8, 26).Value = 0
$ws.Cells.Item(8, 27).Value = 0

# Row 9
$ws.Cells.Item(9, 1).Value = "1-monolinolein, 2tms derivative"
$ws.Cells.Item(9, 2).Value = "1-monolinolein"
$ws.Cells.Item(9, 3).Value = "2,3-dihydroxypropyl (9z,12z)-octadeca-9,12-dienoate"
$ws.Cells.Item(9, 4).Value = "C21H38O4"
$ws.Cells.Item(9, 5).Value = "CCCCCC=CCC=CCCCCCCCC(=O)OCC(CO)O"
$ws.Cells.Item(9, 6).Value = 354.5
$ws.Cells.Item(9, 7).Value = 5.8
$ws.Cells.Item(9, 8).Value = 21
$ws.Cells.Item(9, 9).Value = 38
$ws.Cells.Item(9, 10).Value = 4
$ws.Cells.Item(9, 11).Value = 0.7115119887165021
$ws.Cells.Item(9, 12).Value = 0.1080507757404796
$ws.Cells.Item(9, 13).Value = 0.180524682651622
$ws.Cells.Item(9, 14).Value = 19
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = 2
$ws.Cells.Item(9, 17).Value = 0
$ws.Cells.Item(9, 18).Value = 1
$ws.Cells.Item(9, 19).Value = 0
$ws.Cells.Item(9, 20).Value = 0
$ws.Cells.Item(9, 21).Value = 0.740425952045134
$ws.Cells.Item(9, 22).Value = 0
$ws.Cells.Item(9, 23).Value = 0.09594922425952046
$ws.Cells.Item(9, 24).Value = 0
$ws.Cells.Item(9, 25).Value = 0.1637122708039492
$ws.Cells.Item(9, 26).Value = 0
$ws.Cells.Item(9, 27).Value = 0

# Row 10
$ws.Cells.Item(10, 1).Value = "tetradecanoic acid, deriv"
$ws.Cells.Item(10, 2).Value = "tetradecanoic acid"
$ws.Cells.Item(10, 3).Value = "tetradecanoic acid"
$ws.Cells.Item(10, 4).Value = "C14H28O2"
$ws.Cells.Item(10, 5).Value = "CCCCCCCCCCCCCC(=O)O"
$ws.Cells.Item(10, 6).Value = 228.37
$ws.Cells.Item(10, 7).Value = 5.3
$ws.Cells.Item(10, 8).Value = 14
$ws.Cells.Item(10, 9).Value = 28
$ws.Cells.Item(10, 10).Value = 2
$ws.Cells.Item(10, 11).Value = 0.7363226343214958
$ws.Cells.Item(10, 12).Value = 0.1235889127293427
$ws.Cells.Item(10, 13).Value = 0.1401147261023777
$ws.Cells.Item(10, 14).Value = 13
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = 0
$ws.Cells.Item(10, 19).Value = 0
$ws.Cells.Item(10, 20).Value = 0
$ws.Cells.Item(10, 21).Value = 0.8029031834303979
$ws.Cells.Item(10, 22).Value = 0
$ws.Cells.Item(10, 23).Value = 0
$ws.Cells.Item(10, 24).Value = 0.1971230897228183
$ws.Cells.Item(10, 25).Value = 0
$ws.Cells.Item(10, 26).Value = 0
$ws.Cells.Item(10, 27).Value = 0

# Row 11
$ws.Cells.Item(11, 1).Value = "4-oxopentanoic acid, deriv"
$ws.Cells.Item(11, 2).Value = "4-oxopentanoic acid"
$ws.Cells.Item(11, 3).Value = "4-oxopentanoic acid"
$ws.Cells.Item(11, 4).Value = "C5H8O3"
$ws.Cells.Item(11, 5).Value = "CC(=O)CCC(=O)O"
$ws.Cells.Item(11, 6).Value = 116.11
$ws.Cells.Item(11, 7).Value = -0.5
$ws.Cells.Item(11, 8).Value = 5
$ws.Cells.Item(11, 9).Value = 8
$ws.Cells.Item(11, 10).Value = 3
$ws.Cells.Item(11, 11).Value = 0.5172250452157436
$ws.Cells.Item(11, 12).Value = 0.06945138230987856
$ws.Cells.Item(11, 13).Value = 0.413375247610025
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = 0
$ws.Cells.Item(11, 19).Value = 0
$ws.Cells.Item(11, 20).Value = 1
$ws.Cells.Item(11, 21).Value = 0.1208078546206184
$ws.Cells.Item(11, 22).Value = 0
$ws.Cells.Item(11, 23).Value = 0
$ws.Cells.Item(11, 24).Value = 0.3877099302385669
$ws.Cells.Item(11, 25).Value = 0
$ws.Cells.Item(11, 26).Value = 0
$ws.Cells.Item(11, 27).Value = 0.491533890276462

# Row 12
$ws.Cells.Item(12, 1).Value = "palmitic acid, tms derivative"
$ws.Cells.Item(12, 2).Value = "palmitic acid"
$ws.Cells.Item(12, 3).Value = "hexadecanoic acid"
$ws.Cells.Item(12, 4).Value = "C16H32O2"
$ws.Cells.Item(12, 5).Value = "CCCCCCCCCCCCCCCC(=O)O"
$ws.Cells.Item(12, 6).Value = 256.42
$ws.Cells.Item(12, 7).Value = 6.4
$ws.Cells.Item(12, 8).Value = 16
$ws.Cells.Item(12, 9).Value = 32
$ws.Cells.Item(12, 10).Value = 2
$ws.Cells.Item(12, 11).Value = 0.7494579205990172
$ws.Cells.Item(12, 12).Value = 0.125793619842446
$ws.Cells.Item(12, 13).Value = 0.1247874580765931
$ws.Cells.Item(12, 14).Value = 15
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = 0
$ws.Cells.Item(12, 19).Value = 0
$ws.Cells.Item(12, 20).Value = 0
$ws.Cells.Item(12, 21).Value = 0.8244793697839481
$ws.Cells.Item(12, 22).Value = 0
$ws.Cells.Item(12, 23).Value = 0
$ws.Cells.Item(12, 24).Value = 0.1755596287341081
$ws.Cells.Item(12, 25).Value = 0
$ws.Cells.Item(12, 26).Value = 0
$ws.Cells.Item(12, 27).Value = 0

# Row 13
$ws.Cells.Item(13, 1).Value = "hexadecanoic acid, deriv"
$ws.Cells.Item(13, 2).Value = "hexadecanoic acid"
$ws.Cells.Item(13, 3).Value = "hexadecanoic acid"
$ws.Cells.Item(13, 4).Value = "C16H32O2"
$ws.Cells.Item(13, 5).Value = "CCCCCCCCCCCCCCCC(=O)O"
$ws.Cells.Item(13, 6).Value = 256.42
$ws.Cells.Item(13, 7).Value = 6.4
$ws.Cells.Item(13, 8).Value = 16
$ws.Cells.Item(13, 9).Value = 32
$ws.Cells.Item(13, 10).Value = 2
$ws.Cells.Item(13, 11).Value = 0.7494579205990172
$ws.Cells.Item(13, 12).Value = 0.125793619842446
$ws.Cells.Item(13, 13).Value = 0.1247874580765931
$ws.Cells.Item(13, 14).Value = 15
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 19).Value = 0
$ws.Cells.Item(13, 20).Value = 0
$ws.Cells.Item(13, 21).Value = 0.8244793697839481
$ws.Cells.Item(13, 22).Value = 0
$ws.Cells.Item(13, 23).Value = 0
$ws.Cells.Item(13, 24).Value = 0.1755596287341081
$ws.Cells.Item(13, 25).Value = 0
$ws.Cells.Item(13, 26).Value = 0
$ws.Cells.Item(13, 27).Value = 0

# Row 14
$ws.Cells.Item(14, 1).Value = "2-methoxy-4-methylphenol, deriv"
$ws.Cells.Item(14, 2).Value = "2-methoxy-4-methylphenol"
$ws.Cells.Item(14, 3).Value = "2-methoxy-4-methylphenol"
$ws.Cells.Item(14, 4).Value = "C8H10O2"
$ws.Cells.Item(14, 5).Value = "CC1=CC(=C(C=C1)O)OC"
$ws.Cells.Item(14, 6).Value = 138.16
$ws.Cells.Item(14, 7).Value = 1.3
$ws.Cells.Item(14, 8).Value = 8
$ws.Cells.Item(14, 9).Value = 10
$ws.Cells.Item(14, 10).Value = 2
$ws.Cells.Item(14, 11).Value = 0.6954834973943254
$ws.Cells.Item(14, 12).Value = 0.07295888824551246
$ws.Cells.Item(14, 13).Value = 0.2316010422698321
$ws.Cells.Item(14, 14).Value = 1
$ws.Cells.Item(14, 15).Value = 5
$ws.Cells.Item(14, 16).Value = 1
$ws.Cells.Item(14, 17).Value = 0
$ws.Cells.Item(14, 18).Value = 0
$ws.Cells.Item(14, 19).Value = 1
$ws.Cells.Item(14, 20).Value = 0
$ws.Cells.Item(14, 21).Value = 0.1088231036479444
$ws.Cells.Item(14, 22).Value = 0.4565648523451071
$ws.Cells.Item(14, 23).Value = 0.1230964099594673
$ws.Cells.Item(14, 24).Value = 0
$ws.Cells.Item(14, 25).Value = 0
$ws.Cells.Item(14, 26).Value = 0.3115590619571512
$ws.Cells.Item(14, 27).Value = 0

# Row 15
$ws.Cells.Item(15, 1).Value = "2,4-dimethylphenol, deriv"
$ws.Cells.Item(15, 2).Value = "2,4-dimethylphenol"
$ws.Cells.Item(15, 3).Value = "2,4-dimethylphenol"
$ws.Cells.Item(15, 4).Value = "C8H10O"
$ws.Cells.Item(15, 5).Value = "CC1=CC(=C(C=C1)O)C"
$ws.Cells.Item(15, 6).Value = 122.16
$ws.Cells.Item(15, 7).Value = 2.3
$ws.Cells.Item(15, 8).Value = 8
$ws.Cells.Item(15, 9).Value = 10
$ws.Cells.Item(15, 10).Value = 1
$ws.Cells.Item(15, 11).Value = 0.7865749836280288
$ws.Cells.Item(15, 12).Value = 0.0825147347740668
$ws.Cells.Item(15, 13).Value = 0.1309675834970531
$ws.Cells.Item(15, 14).Value = 2
$ws.Cells.Item(15, 15).Value = 6
$ws.Cells.Item(15, 16).Value = 1
$ws.Cells.Item(15, 17).Value = 0
$ws.Cells.Item(15, 18).Value = 0
$ws.Cells.Item(15, 19).Value = 0
$ws.Cells.Item(15, 20).Value = 0
$ws.Cells.Item(15, 21).Value = 0.2461525867714473
$ws.Cells.Item(15, 22).Value = 0.6146856581532416
$ws.Cells.Item(15, 23).Value = 0.1392190569744597
$ws.Cells.Item(15, 24).Value = 0
$ws.Cells.Item(15, 25).Value = 0
$ws.Cells.Item(15, 26).Value = 0
$ws.Cells.Item(15, 27).Value = 0

# Row 16
$ws.Cells.Item(16, 1).Value = "palmitelaidic acid, tms derivative"
$ws.Cells.Item(16, 2).Value = "palmitelaidic acid"
$ws.Cells.Item(16, 3).Value = "(e)-hexadec-9-enoic acid"
$ws.Cells.Item(16, 4).Value = "C16H30O2"
$ws.Cells.Item(16, 5).Value = "CCCCCCC=CCCCCCCCC(=O)O"
$ws.Cells.Item(16, 6).Value = 254.41
$ws.Cells.Item(16, 7).Value = 6.4
$ws.Cells.Item(16, 8).Value = 16
$ws.Cells.Item(16, 9).Value = 30
$ws.Cells.Item(16, 10).Value = 2
$ws.Cells.Item(16, 11).Value = 0.7553791124562713
$ws.Cells.Item(16, 12).Value = 0.1188632522306513
$ws.Cells.Item(16, 13).Value = 0.1257733579654888
$ws.Cells.Item(16, 14).Value = 15
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = 0
$ws.Cells.Item(16, 19).Value = 0
$ws.Cells.Item(16, 20).Value = 0
$ws.Cells.Item(16, 21).Value = 0.8230690617507173
$ws.Cells.Item(16, 22).Value = 0
$ws.Cells.Item(16, 23).Value = 0
$ws.Cells.Item(16, 24).Value = 0.1769466609016941
$ws.Cells.Item(16, 25).Value = 0
$ws.Cells.Item(16, 26).Value = 0
$ws.Cells.Item(16, 27).Value = 0

# Row 17
$ws.Cells.Item(17, 1).Value = "(9z,12z)-octadeca-9,12-dienoic acid, deriv"
$ws.Cells.Item(17, 2).Value = "(9z,12z)-octadeca-9,12-dienoic acid"
$ws.Cells.Item(17, 3).Value = "(9z,12z)-octadeca-9,12-dienoic acid"
$ws.Cells.Item(17, 4).Value = "C18H32O2"
$ws.Cells.Item(17, 5).Value = "CCCCCC=CCC=CCCCCCCCC(=O)O"
$ws.Cells.Item(17, 6).Value = 280.4
$ws.Cells.Item(17, 7).Value = 6.8
$ws.Cells.Item(17, 8).Value = 18
$ws.Cells.Item(17, 9).Value = 32
$ws.Cells.Item(17, 10).Value = 2
$ws.Cells.Item(17, 11).Value = 0.7710342368045648
$ws.Cells.Item(17, 12).Value = 0.1150356633380885
$ws.Cells.Item(17, 13).Value = 0.1141155492154066
$ws.Cells.Item(17, 14).Value = 17
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = 0
$ws.Cells.Item(17, 19).Value = 0
$ws.Cells.Item(17, 20).Value = 0
$ws.Cells.Item(17, 21).Value = 0.8396398002853066
$ws.Cells.Item(17, 22).Value = 0
$ws.Cells.Item(17, 23).Value = 0
$ws.Cells.Item(17, 24).Value = 0.1605456490727532
$ws.Cells.Item(17, 25).Value = 0
$ws.Cells.Item(17, 26).Value = 0
$ws.Cells.Item(17, 27).Value = 0

# Row 18
$ws.Cells.Item(18, 1).Value = "myristic acid, tms derivative"
$ws.Cells.Item(18, 2).Value = "myristic acid"
$ws.Cells.Item(18, 3).Value = "tetradecanoic acid"
$ws.Cells.Item(18, 4).Value = "C14H28O2"
$ws.Cells.Item(18, 5).Value = "CCCCCCCCCCCCCC(=O)O"
$ws.Cells.Item(18, 6).Value = 228.37
$ws.Cells.Item(18, 7).Value = 5.3
$ws.Cells.Item(18, 8).Value = 14
$ws.Cells.Item(18, 9).Value = 28
$ws.Cells.Item(18, 10).Value = 2
$ws.Cells.Item(18, 11).Value = 0.7363226343214958
$ws.Cells.Item(18, 12).Value = 0.1235889127293427
$ws.Cells.Item(18, 13).Value = 0.1401147261023777
$ws.Cells.Item(18, 14).Value = 13
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 1
$ws.Cells.Item(18, 18).Value = 0
$ws.Cells.Item(18, 19).Value = 0
$ws.Cells.Item(18, 20).Value = 0
$ws.Cells.Item(18, 21).Value = 0.8029031834303979
$ws.Cells.Item(18, 22).Value = 0
$ws.Cells.Item(18, 23).Value = 0
$ws.Cells.Item(18, 24).Value = 0.1971230897228183
$ws.Cells.Item(18, 25).Value = 0
$ws.Cells.Item(18, 26).Value = 0
$ws.Cells.Item(18, 27).Value = 0

# Row 19
$ws.Cells.Item(19, 1).Value = "phenol, deriv"
$ws.Cells.Item(19, 2).Value = "phenol"
$ws.Cells.Item(19, 3).Value = "phenol"
$ws.Cells.Item(19, 4).Value = "C6H6O"
$ws.Cells.Item(19, 5).Value = "C1=CC=C(C=C1)O"
$ws.Cells.Item(19, 6).Value = 94.11
$ws.Cells.Item(19, 7).Value = 1.5
$ws.Cells.Item(19, 8).Value = 6
$ws.Cells.Item(19, 9).Value = 6
$ws.Cells.Item(19, 10).Value = 1
$ws.Cells.Item(19, 11).Value = 0.765763468281798
$ws.Cells.Item(19, 12).Value = 0.06426522154925088
$ws.Cells.Item(19, 13).Value = 0.1700031877590054
$ws.Cells.Item(19, 14).Value = 0
$ws.Cells.Item(19, 15).Value = 6
$ws.Cells.Item(19, 16).Value = 1
$ws.Cells.Item(19, 17).Value = 0
$ws.Cells.Item(19, 18).Value = 0
$ws.Cells.Item(19, 19).Value = 0
$ws.Cells.Item(19, 20).Value = 0
$ws.Cells.Item(19, 21).Value = 0
$ws.Cells.Item(19, 22).Value = 0.8193178195728402
$ws.Cells.Item(19, 23).Value = 0.1807140580172139
$ws.Cells.Item(19, 24).Value = 0
$ws.Cells.Item(19, 25).Value = 0
$ws.Cells.Item(19, 26).Value = 0
$ws.Cells.Item(19, 27).Value = 0

# Row 20
$ws.Cells.Item(20, 1).Value = "glycerol monostearate, 2tms derivative"
$ws.Cells.Item(20, 2).Value = "glycerol monostearate"
$ws.Cells.Item(20, 3).Value = "2,3-dihydroxypropyl octadecanoate"
$ws.Cells.Item(20, 4).Value = "C21H42O4"
$ws.Cells.Item(20, 5).Value = "CCCCCCCCCCCCCCCCCC(=O)OCC(CO)O"
$ws.Cells.Item(20, 6).Value = 358.6
$ws.Cells.Item(20, 7).Value = 7.4
$ws.Cells.Item(20, 8).Value = 21
$ws.Cells.Item(20, 9).Value = 42
$ws.Cells.Item(20, 10).Value = 4
$ws.Cells.Item(20, 11).Value = 0.7033770217512548
$ws.Cells.Item(20, 12).Value = 0.1180591187953151
$ws.Cells.Item(20, 13).Value = 0.1784606804238706
$ws.Cells.Item(20, 14).Value = 19
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = 2
$ws.Cells.Item(20, 17).Value = 0
$ws.Cells.Item(20, 18).Value = 1
$ws.Cells.Item(20, 19).Value = 0
$ws.Cells.Item(20, 20).Value = 0
$ws.Cells.Item(20, 21).Value = 0.7432041271611823
$ws.Cells.Item(20, 22).Value = 0
$ws.Cells.Item(20, 23).Value = 0.09485220301171221
$ws.Cells.Item(20, 24).Value = 0
$ws.Cells.Item(20, 25).Value = 0.161840490797546
$ws.Cells.Item(20, 26).Value = 0
$ws.Cells.Item(20, 27).Value = 0

# Row 21
$ws.Cells.Item(21, 1).Value = "oleic acid, tms derivative"
$ws.Cells.Item(21, 2).Value = "oleic acid"
$ws.Cells.Item(21, 3).Value = "(z)-octadec-9-enoic acid"
$ws.Cells.Item(21, 4).Value = "C18H34O2"
$ws.Cells.Item(21, 5).Value = "CCCCCCCCC=CCCCCCCCC(=O)O"
$ws.Cells.Item(21, 6).Value = 282.5
$ws.Cells.Item(21, 7).Value = 6.5
$ws.Cells.Item(21, 8).Value = 18
$ws.Cells.Item(21, 9).Value = 34
$ws.Cells.Item(21, 10).Value = 2
$ws.Cells.Item(21, 11).Value = 0.7653026548672566
$ws.Cells.Item(21, 12).Value = 0.121316814159292
$ws.Cells.Item(21, 13).Value = 0.1132672566371681
$ws.Cells.Item(21, 14).Value = 17
$ws.Cells.Item(21, 15).Value = 0
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = 0
$ws.Cells.Item(21, 19).Value = 0
$ws.Cells.Item(21, 20).Value = 0
$ws.Cells.Item(21, 21).Value = 0.8405345132743363
$ws.Cells.Item(21, 22).Value = 0
$ws.Cells.Item(21, 23).Value = 0
$ws.Cells.Item(21, 24).Value = 0.1593522123893805
$ws.Cells.Item(21, 25).Value = 0
$ws.Cells.Item(21, 26).Value = 0
$ws.Cells.Item(21, 27).Value = 0

# Row 22
$ws.Cells.Item(22, 1).Value = "benzoic acid, deriv"
$ws.Cells.Item(22, 2).Value = "benzoic acid"
$ws.Cells.Item(22, 3).Value = "benzoic acid"
$ws.Cells.Item(22, 4).Value = "C7H6O2"
$ws.Cells.Item(22, 5).Value = "C1=CC=C(C=C1)C(=O)O"
$ws.Cells.Item(22, 6).Value = 122.12
$ws.Cells.Item(22, 7).Value = 1.9
$ws.Cells.Item(22, 8).Value = 7
$ws.Cells.Item(22, 9).Value = 6
$ws.Cells.Item(22, 10).Value = 2
$ws.Cells.Item(22, 11).Value = 0.6884785456927611
$ws.Cells.Item(22, 12).Value = 0.04952505732066819
$ws.Cells.Item(22, 13).Value = 0.2620209629872257
$ws.Cells.Item(22, 14).Value = 0
$ws.Cells.Item(22, 15).Value = 6
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = 0
$ws.Cells.Item(22, 19).Value = 0
$ws.Cells.Item(22, 20).Value = 0
$ws.Cells.Item(22, 21).Value = 0
$ws.Cells.Item(22, 22).Value = 0.6313953488372093
$ws.Cells.Item(22, 23).Value = 0
$ws.Cells.Item(22, 24).Value = 0.3686292171634458
$ws.Cells.Item(22, 25).Value = 0
$ws.Cells.Item(22, 26).Value = 0
$ws.Cells.Item(22, 27).Value = 0

